$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a truly-blank-but-present text cell (matches the source
# workbook's convention of emitting an empty inline/shared string rather
# than omitting the cell entirely). Entering "'" forces Excel to commit a
# Text-typed empty string, then resetting the Style back to Normal drops
# the quote-prefix formatting so the cell carries no extra style.
function Set-BlankTextCell($addr) {
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("A9").Value = "gemini-1.5-pro"
Set-BlankTextCell "B9"
$ws.Range("C9").Value = "0.21 ± 0.71"
$ws.Range("D9").Value = "0.45 ± 0.76"
Set-BlankTextCell "E9"
Set-BlankTextCell "F9"
Set-BlankTextCell "G9"
Set-BlankTextCell "H9"
Set-BlankTextCell "I9"
$ws.Range("J9").Value = "0.79 ± 0.18"
$ws.Range("K9").Value = "0.81 ± 0.19"
$ws.Range("L9").Value = "0.8 ± 0.19"
$ws.Range("M9").Value = "0.85 ± 0.2"
$ws.Range("N9").Value = "0.95 ± 0.22"
Set-BlankTextCell "O9"
$ws.Range("P9").Value = "0.62 ± 0.18"
$ws.Range("Q9").Value = "2.89 ± 3.99"
$ws.Range("R9").Value = "0.051 ± 0.00"
$ws.Range("S9").Value = "0.9 ± 0.21"
$ws.Range("T9").Value = "0.87 ± 0.27"
$ws.Range("U9").Value = "3.15 ± 1.49"
$ws.Range("V9").Value = "0.88 ± 0.27"
$ws.Range("W9").Value = "0.89 ± 0.21"
$ws.Range("X9").Value = "1.34 ± 0.41"
